$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$url = "http://www.digikey.com/scripts/DkSearch/dksus.dll?Detail&itemSeq=178538991&uq=635760440624849824"

# Row 6: add a new capacitor line (445-8919-1-ND, linked to Digikey) with its note
$ws.Hyperlinks.Add($ws.Range("C6"), $url, "", "", $url)
$ws.Range("C6").Value = "445-8919-1-ND"
$ws.Range("E6").Value = "10 uF decoupling"

# Row 11: add another capacitor line (490-6125-1-ND), matching the P/N style used elsewhere
$ws.Range("A11").Value = "Capacitor"
$ws.Range("C11").Value = "490-6125-1-ND"
$ws.Range("C9").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("D11").Value = 1
$ws.Application.CutCopyMode = $false

# Update the active selection to match the edited area
$ws.Range("E11").Select() | Out-Null
